# CIERRE 6 MAY 22
# - Update the incentive-payment month/year caption on "VALES DE INSENTIVOS"
#   from "MARZO 2022" to " ABRIL  2022".
# - Move the selected/active tab from "ARQUITECTO" to "VALES DE INSENTIVOS".
# (The TODAY() cells on both sheets re-cache automatically on recalc.)

$wb = $excel.ActiveWorkbook

$wsVales = $wb.Worksheets.Item("VALES DE INSENTIVOS")
$wsVales.Range("A4").Value = "PAGO DE INCENTIVO DEL MES DE  ABRIL   2022"

# Activating this sheet makes it the workbook's active tab (activeTab on the
# workbookView) and moves tabSelected from "ARQUITECTO" to this sheet.
$wsVales.Activate()
